$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '69.463.73'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.68%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.689.59'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.18%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '686.13'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '160.06'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("E9").Value = '  -0.62%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '7.09'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  -2.79%  '
$ws.Range("E12").Value = '  -0.36%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.315.15'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.41%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '32.54'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -2.19%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.705.29'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.01%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '69.430.81'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").Value = '  +1.73%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '15.85'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.39%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '6.42'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -2.06%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '471.23'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.56%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '10.02'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("E22").Value = '  -1.50%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '79.75'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +1.12%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '3.836.41'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -2.20%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '11.02'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -4.61%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.26'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("E31").Value = '  -4.55%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '6.58'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -1.57%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.23%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '26.95'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.95%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.663.12'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.84%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.160'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.53%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '8.20'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -3.12%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '6.16'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("E40").Value = '  +2.11%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0899'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -3.26%  '
$ws.Range("E42").Value = '  +0.07%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.943'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("E44").Value = '  +2.68%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '47.56'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -1.65%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.000283'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +3.77%  '
$ws.Range("E47").Value = '  +7.44%  '
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("E49").Value = '  -0.47%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '27.94'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.01%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '7.77'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.56%  '
